# Re-applies the latest cryptocurrency market snapshot onto the
# "cryptos" sheet: updated Price (D) / Volume(1h) (E) figures, plus
# three pairs of rows whose rank swapped places (B/C/D/E all move).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address -> new text. A leading "'" forces the cell to
# stay plain text (many prices look numeric, e.g. "1.11" or "96.816.73",
# and would otherwise be auto-coerced to a Double by Excel); re-applying
# the "Normal" style afterwards clears the resulting quote-prefix marker
# so the cell format matches the original unstyled cells exactly.
$updates = [ordered]@{
    "D2" = "96.816.73"
    "E2" = "  +0.64%  "
    "D3" = "3.711.21"
    "E3" = "  +3.93%  "
    "E4" = "  -0.02%  "
    "D5" = "243.43"
    "E5" = "  +1.18%  "
    "D6" = "1.91"
    "E6" = "  +19.91%  "
    "D7" = "660.12"
    "E7" = "  +0.96%  "
    "D8" = "0.428"
    "E8" = "  +5.85%  "
    "D9" = "1.11"
    "E9" = "  +6.40%  "
    "D10" = "1.00"
    "E10" = "  -0.03%  "
    "D11" = "3.709.97"
    "E11" = "  +3.93%  "
    "D12" = "44.94"
    "E12" = "  +4.45%  "
    "E13" = "  +1.58%  "
    "D14" = "6.54"
    "E14" = "  +2.44%  "
    "D15" = "4.402.18"
    "E15" = "  +3.92%  "
    "D16" = "96.821.97"
    "E16" = "  +0.73%  "
    "E17" = "  +1.25%  "
    "D18" = "3.713.91"
    "E18" = "  +4.17%  "
    "B19" = "Uniswap"
    "C19" = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
    "D19" = "13.30"
    "E19" = "  +6.01%  "
    "B20" = "Polkadot"
    "C20" = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
    "D20" = "7.78"
    "E20" = "  +0.44%  "
    "E21" = "  +5.08%  "
    "D22" = "0.536"
    "E22" = "  +5.79%  "
    "D23" = "3.46"
    "E23" = "  +1.37%  "
    "D24" = "514.14"
    "E24" = "  +1.30%  "
    "D25" = "0.0000213"
    "E25" = "  +6.86%  "
    "D26" = "6.89"
    "E26" = "  +0.63%  "
    "E27" = "  +5.72%  "
    "D28" = "13.14"
    "E28" = "  +4.18%  "
    "D29" = "0.171"
    "E29" = "  +14.01%  "
    "D30" = "3.06"
    "E30" = "  +2.80%  "
    "D31" = "12.05"
    "E31" = "  +5.83%  "
    "E32" = "  +0.24%  "
    "E33" = "  +2.18%  "
    "D34" = "33.47"
    "E34" = "  +6.43%  "
    "D35" = "1.00"
    "E35" = "  -0.06%  "
    "B36" = "PolygonEcosystemToken"
    "C36" = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
    "D36" = "0.595"
    "E36" = "  +5.54%  "
    "B37" = "Fetch.AI"
    "C37" = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
    "D37" = "1.72"
    "E37" = "  +6.37%  "
    "D38" = "615.58"
    "E38" = "  -1.08%  "
    "D39" = "8.74"
    "E39" = "  +0.00%  "
    "D40" = "42.86"
    "E40" = "  +26.20%  "
    "E41" = "  +5.84%  "
    "D42" = "0.973"
    "E42" = "  +7.70%  "
    "E43" = "  +7.05%  "
    "E44" = "  +0.01%  "
    "D45" = "6.13"
    "E45" = "  +8.07%  "
    "E46" = "  +5.34%  "
    "B47" = "Algorand"
    "C47" = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
    "D47" = "0.419"
    "E47" = "  +24.76%  "
    "B48" = "Stacks"
    "C48" = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
    "D48" = "2.31"
    "E48" = "  +0.89%  "
    "D49" = "23.59"
    "E49" = "  +0.41%  "
    "D50" = "8.59"
    "E50" = "  +5.35%  "
    "D51" = "54.53"
    "E51" = "  +2.91%  "
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $updates[$addr]
    $cell.Style = "Normal"
}
